$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ether" column (column B) entirely, shifting C:E left to B:D
$ws.Range("B1").EntireColumn.Delete()

# Update header row (now only 4 columns: amines, aldehyde, aromatic, cycle)
$ws.Range("A1").Value = "amines"
$ws.Range("B1").Value = "aldehyde"
$ws.Range("C1").Value = "aromatic"
$ws.Range("D1").Value = "cycle"

# Update the numeric data values
$ws.Range("A2").Value = 1.683745353784681
$ws.Range("B2").Value = 0.240198009795508
$ws.Range("C2").Value = 1.021669980951205
$ws.Range("D2").Value = 0.9137032076671433

$ws.Range("A3").Value = 3.460053772526015
$ws.Range("B3").Value = 0.9596744520597827
$ws.Range("C3").Value = 0.6951225421178187
$ws.Range("D3").Value = 1.618744173530442

$ws.Range("A4").Value = 0.8181676205552553
$ws.Range("B4").Value = 1.762780576214276
$ws.Range("C4").Value = 0.9355905554710553
$ws.Range("D4").Value = 3.676130547679126
